# Update "想去人数" (wanted-to-go count) values in F column for both the
# "展览" sheet and the consolidated "全部类型" sheet, matching the rows that
# refer to the same events.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 508
$ws1.Range("F7").Value = 1018
$ws1.Range("F11").Value = 216
$ws1.Range("F13").Value = 182
$ws1.Range("F14").Value = 168

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 508
$ws4.Range("F8").Value = 1018
$ws4.Range("F12").Value = 216
$ws4.Range("F14").Value = 182
$ws4.Range("F15").Value = 168
